$d = $word.ActiveDocument

# --- Step 1: Fill in the placeholder single-character cell values in the
#             5 pre-existing Risk tables with the full report text. ---
$t = $d.Tables.Item(1)
$t.Cell(1, 2).Range.Text = "The Linux Exploit Suggester identified several potential vulnerabilities (CVEs) in the system's kernel and sudo configuration.  While the likelihood of successful exploitation depends on various factors, the presence of these vulnerabilities poses a significant risk."
$t.Cell(2, 2).Range.Text = "High"
$t.Cell(3, 2).Range.Text = "Very High"
$t.Cell(4, 2).Range.Text = "Attackers could potentially exploit vulnerabilities in the system's kernel or sudo configuration to gain root privileges."
$t.Cell(5, 2).Range.Text = "Update your system's kernel and sudo to the latest versions.  Ensure that all security patches are applied.  Contact your IT administrator or security professional for guidance on addressing these vulnerabilities."

$t = $d.Tables.Item(2)
$t.Cell(1, 2).Range.Text = "The Linpeas output revealed the presence of network discovery and port scanning tools (fping, bash, nc, nmap).  An attacker could use these tools to map the system's network, identify open ports, and potentially exploit vulnerabilities."
$t.Cell(2, 2).Range.Text = "Medium"
$t.Cell(3, 2).Range.Text = "High"
$t.Cell(4, 2).Range.Text = "Attackers could use these tools to map the system's network, identify open ports, and potentially exploit vulnerabilities."
$t.Cell(5, 2).Range.Text = "Disable or remove any unnecessary network discovery and port scanning tools.  Restrict network access to the system as much as possible.  Consult your IT administrator for assistance in securing your network."

$t = $d.Tables.Item(3)
$t.Cell(1, 2).Range.Text = "The Linpeas output shows several writable configuration files, including systemd service files.  An attacker could modify these files to compromise the system."
$t.Cell(2, 2).Range.Text = "High"
$t.Cell(3, 2).Range.Text = "High"
$t.Cell(4, 2).Range.Text = "An attacker could modify these files to compromise the system."
$t.Cell(5, 2).Range.Text = "Make all critical configuration files read-only.  Implement appropriate access control measures to prevent unauthorized modification of these files.  Work with your IT administrator to establish secure configuration management practices."

$t = $d.Tables.Item(4)
$t.Cell(1, 2).Range.Text = "The Linpeas output shows several software packages that have known vulnerabilities or are outdated.  These vulnerabilities could be exploited by attackers."
$t.Cell(2, 2).Range.Text = "Medium"
$t.Cell(3, 2).Range.Text = "High"
$t.Cell(4, 2).Range.Text = "These vulnerabilities could be exploited by attackers."
$t.Cell(5, 2).Range.Text = "Update all software packages to their latest versions.  Regularly check for and install security updates.  Use a vulnerability scanner to identify and address any known vulnerabilities."

$t = $d.Tables.Item(5)
$t.Cell(1, 2).Range.Text = "The Linpeas output shows that several ports are open.  These open ports could be exploited by attackers to gain unauthorized access to the system."
$t.Cell(2, 2).Range.Text = "Medium"
$t.Cell(3, 2).Range.Text = "High"
$t.Cell(4, 2).Range.Text = "These open ports could be exploited by attackers to gain unauthorized access to the system."
$t.Cell(5, 2).Range.Text = "Close any unnecessary ports.  Use a firewall to restrict network access to only necessary ports.  Consult your IT administrator for assistance in configuring your firewall."

# --- Step 2: Append 3 new Risk tables (same 2-col/5-row shape) after the
#             last existing table, each preceded by its own blank paragraph
#             exactly like the existing tables are separated. ---
$lastTable = $d.Tables.Item($d.Tables.Count)
$insertRange = $d.Range($lastTable.Range.End, $lastTable.Range.End)

$newTablesXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:tbl><w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="4320"/><w:gridCol w:w="4320"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Statement</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>The Linpeas scan revealed the presence of sensitive data in various files, such as password hashes and API keys (if the '-r' parameter had been used).  This data could be exploited by attackers.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Likelihood</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>High</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Impact</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>High</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Impact of Risk on system</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>This data could be exploited by attackers.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>What to do</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Securely store all sensitive data, such as password hashes and API keys.  Use strong passwords and encryption to protect this data.  Consult your IT administrator for guidance on secure data handling practices.</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="4320"/><w:gridCol w:w="4320"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Statement</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Multiple vulnerabilities were identified in the system, including the presence of writable files in critical directories (/etc/passwd), SUID/SGID binaries with potential for privilege escalation, and the availability of network discovery and port scanning tools.  These vulnerabilities, combined with the lack of several security protections, create a high likelihood of successful exploitation.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Likelihood</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>High</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Impact</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Very High</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Impact of Risk on system</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Attackers could potentially exploit vulnerabilities in the system's kernel or sudo configuration to gain root privileges.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>What to do</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Because you are already root and you have identified several risks, you should immediately contact your IT administrator or security professional.  They can help you secure your system and prevent unauthorized access.  Do not attempt to fix these issues yourself unless you are an experienced IT professional.</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="4320"/><w:gridCol w:w="4320"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Statement</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>The Linpeas output shows environment variables containing sensitive information such as SSH_AGENT_PID and XAUTHORITY.  While not directly exploitable, this information could aid an attacker in further compromising the system.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Likelihood</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Medium</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Risk Impact</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Medium</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Impact of Risk on system</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>This information could aid an attacker in further compromising the system.</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>What to do</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="4320"/></w:tcPr><w:p><w:r><w:t>Review the environment variables listed in the Linpeas output.  If any contain sensitive information, remove or change them immediately.  Consult your IT administrator for assistance if needed.</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$insertRange.InsertXML($newTablesXml)

